# Automatic update of files.
#
# Applies the diff to "A 1480-2026 artfynd.xlsx": a handful of observation
# rows had their species/location details rotated between rows (same
# survey, same date/observer/site, different taxon pulled into each row),
# plus several unrelated standalone "Taxonsorteringsordning" (column B)
# bumps of +3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 7 <-> 8 : full species-record swap (Brunpudrad nållav <-> Tretåig
# hackspett). Row 7 gains the bird-specific K/L/M/N + AC cells; row 8
# loses them.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 131066771
$ws.Range("B7").Value = 57884
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "äldre spår"
$ws.Range("Q7").Value = 425297
$ws.Range("R7").Value = 6712214
$ws.Range("AC7").Value = "Ringhack på gran"

$ws.Range("A8").Value = 131066788
$ws.Range("B8").Value = 83217
$ws.Range("E8").Value = 308
$ws.Range("F8").Value = "Brunpudrad nållav"
$ws.Range("G8").Value = "Chaenotheca gracillima"
$ws.Range("H8").Value = "(Vain.) Tibell"
$ws.Range("K8:N8").ClearContents()
$ws.Range("Q8").Value = 425211
$ws.Range("R8").Value = 6712276
$ws.Range("AC8").ClearContents()

# ---------------------------------------------------------------------
# Standalone Taxonsorteringsordning (+3) bumps - only column B changes.
# ---------------------------------------------------------------------
$ws.Range("B11").Value = 91776
$ws.Range("B15").Value = 91776
$ws.Range("B16").Value = 91827
$ws.Range("B17").Value = 91813

# ---------------------------------------------------------------------
# Rows 21 -> 22 -> 23 -> 21 : three-way rotation of species records
# (Ulltickeporing / Kortskaftad ärgspik / Tretåig hackspett), each
# keeping its own site coordinates. Row 23's Taxonsorteringsordning also
# independently bumps by +3 (92181 -> 92184) on top of the rotation.
# ---------------------------------------------------------------------
$ws.Range("A21").Value = 131066778
$ws.Range("B21").Value = 81230
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 1049
$ws.Range("F21").Value = "Kortskaftad ärgspik"
$ws.Range("G21").Value = "Microcalicium ahlneri"
$ws.Range("H21").Value = "Tibell"
$ws.Range("Q21").Value = 425336
$ws.Range("R21").Value = 6712202

$ws.Range("A22").Value = 131066774
$ws.Range("B22").Value = 57884
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("M22").Value = "äldre spår"
$ws.Range("Q22").Value = 425250
$ws.Range("R22").Value = 6712265
$ws.Range("AC22").Value = "Ringhack på gran"

$ws.Range("A23").Value = 131066766
$ws.Range("B23").Value = 92184
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 2062
$ws.Range("F23").Value = "Ulltickeporing"
$ws.Range("G23").Value = "Skeletocutis brevispora"
$ws.Range("H23").Value = "Niemelä"
$ws.Range("K23:N23").ClearContents()
$ws.Range("Q23").Value = 425069
$ws.Range("R23").Value = 6712285
$ws.Range("AC23").ClearContents()

# ---------------------------------------------------------------------
# More standalone Taxonsorteringsordning (+3) bumps.
# ---------------------------------------------------------------------
$ws.Range("B25").Value = 91776
$ws.Range("B26").Value = 91776
$ws.Range("B29").Value = 92232

# ---------------------------------------------------------------------
# Rows 30 <-> 31 : full species-record swap (Skrovellav <-> Tretåig
# hackspett), same pattern as rows 7/8.
# ---------------------------------------------------------------------
$ws.Range("A30").Value = 131066772
$ws.Range("B30").Value = 57884
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "äldre spår"
$ws.Range("Q30").Value = 425301
$ws.Range("R30").Value = 6712219
$ws.Range("AC30").Value = "Ringhack på gran"

$ws.Range("A31").Value = 131066776
$ws.Range("B31").Value = 80351
$ws.Range("E31").Value = 2081
$ws.Range("F31").Value = "Skrovellav"
$ws.Range("G31").Value = "Lobaria scrobiculata"
$ws.Range("H31").Value = "(Scop.) DC."
$ws.Range("K31:N31").ClearContents()
$ws.Range("Q31").Value = 425069
$ws.Range("R31").Value = 6712285
$ws.Range("AC31").ClearContents()

# ---------------------------------------------------------------------
# Final standalone Taxonsorteringsordning (+3) bump.
# ---------------------------------------------------------------------
$ws.Range("B32").Value = 91813
